$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Execute flag for the PROD/orion row to "NO"
$ws.Range("C4").Value = "NO"

# Update the Execute flag for the PROD/gesco row to "YES"
$ws.Range("C9").Value = "YES"

# Move the active cell selection to C9
$ws.Range("C9").Select()
